# Update "想去人数" (interest counts) figures in the F column across the
# four worksheets, matching the latest data refresh (gh-pages output
# generated at 456a3b4).

$wb = $excel.ActiveWorkbook

# 展览 (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 1075
$ws1.Range("F5").Value = 2477
$ws1.Range("F7").Value = 677
$ws1.Range("F11").Value = 683
$ws1.Range("F14").Value = 1453
$ws1.Range("F18").Value = 252

# 演出 (Performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F16").Value = 7

# 本地生活 (Local Life)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F4").Value = 2005

# 全部类型 (All types - combined listing)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 2005
$ws4.Range("F12").Value = 1075
$ws4.Range("F16").Value = 2477
$ws4.Range("F22").Value = 677
$ws4.Range("F27").Value = 683
$ws4.Range("F31").Value = 1453
$ws4.Range("F34").Value = 7
$ws4.Range("F43").Value = 252
